$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = 45854.65412042624

$data = @(
    ,@("0ce5dd49", "Generating new code snippets|Explaining code or concepts|Writing unit tests")
    ,@("2c1001cb", "Generating new code snippets|Debugging existing code|Explaining code or concepts|Writing unit tests")
    ,@("37cc37bf", "I did not choose “Coding / programming help”")
    ,@("43faa0b9", "Generating new code snippets|Explaining code or concepts|Writing unit tests")
    ,@("4abe3e88", "I did not choose “Coding / programming help”")
    ,@("50164f59", "Generating new code snippets|Debugging existing code|Explaining code or concepts|Writing unit tests")
    ,@("5cf70f79", "Generating new code snippets|Debugging existing code|Explaining code or concepts|Writing unit tests")
    ,@("5da96769", "I did not choose “Coding / programming help”")
    ,@("6ca3e2f6", "Explaining code or concepts")
    ,@("790a4fcb", "Generating new code snippets|Explaining code or concepts|Writing unit tests")
    ,@("802cc63a", "Generating new code snippets|Debugging existing code|Explaining code or concepts|Converting code between languages|Writing unit tests")
    ,@("85c3ea4d", "Generating new code snippets|Debugging existing code|Explaining code or concepts|Converting code between languages")
    ,@("942dfafb", "Generating new code snippets|Debugging existing code|Explaining code or concepts|Writing unit tests")
    ,@("9bc6ba8c", "I did not choose “Coding / programming help”")
    ,@("a2d65af2", "Explaining code or concepts")
    ,@("a46f1771", "I did not choose “Coding / programming help”")
    ,@("ad58f9da", "I did not choose “Coding / programming help”")
    ,@("c7d9a301", "Generating new code snippets|Debugging existing code|Explaining code or concepts|Converting code between languages")
    ,@("ce8732ff", "I did not choose “Coding / programming help”")
    ,@("d6f1d567", "Generating new code snippets|Debugging existing code|Explaining code or concepts|Converting code between languages")
    ,@("da9326c9", "Generating new code snippets|Explaining code or concepts")
    ,@("e09ca7bf", "Generating new code snippets|Debugging existing code|Explaining code or concepts|Writing unit tests")
    ,@("ef53a641", "Generating new code snippets|Explaining code or concepts|Writing unit tests")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = "q09_code_subtasks"
    $ws.Cells.Item($row, 4).Value = $timestamp
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $row++
}
